$d = $word.ActiveDocument

# --- 1) "Com base na Tabela..." block (4 paragraphs) -> single merged paragraph with new text ---
# Locate the paragraph that starts the block by its distinctive text.
$startPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Com base na Tabela e os itens analisados*") {
        $startPara = $i
        break
    }
}

if ($startPara -ne $null) {
    # Merge the 4 paragraphs (this one + the next 3) into one by deleting the
    # paragraph-mark characters that separate them.
    for ($k = 0; $k -lt 3; $k++) {
        $p = $d.Paragraphs.Item($startPara)
        $markRange = $d.Range($p.Range.End - 1, $p.Range.End)
        $markRange.Delete()
    }
    $merged = $d.Paragraphs.Item($startPara)
    $textRange = $d.Range($merged.Range.Start, $merged.Range.End - 1)
    $textRange.Text = "Com base na minha experiência pessoal de uso, percebi muita qualidade no material do aparelho, excelente durabilidade da bateria com uso moderado e intenso com vários processos ao mesmo temo, com um processador com mais potência foi notável que não obtive nenhum erro de travamento, e adicionando as qualidades da Câmera frontal e as câmeras traseira com excelente qualidade. Um produto com excelentes pontos positivos."
}

# --- 2) "Imagem 1: " -> "Demonstração de qualidade de Imagem da " (before "Câmera traseira") ---
$d.Content.Find.Execute("Imagem 1: ", $true, $false, $false, $false, $false, $true, 1, $false, "Demonstração de qualidade de Imagem da ", 2) | Out-Null

# --- 3) Remove the standalone highlighted "Foto:" run (leaving the paragraph empty) ---
$d.Content.Find.Execute("Foto:", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- 4) "Produto: Foto Xiaomi Readmi Note 11" -> "Produto: Foto ilustração  Xiaomi Readmi Note 11" ---
$d.Content.Find.Execute(": Foto Xiaomi ", $true, $false, $false, $false, $false, $true, 1, $false, ": Foto ilustração  Xiaomi ", 2) | Out-Null

Write-Host "done"
